$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'261.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'0.86%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('E3').Value = "'0.98%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'4.704"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'0.59%"
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.06197"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'2.96%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'6.727"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'0.93%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.8511"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'-1.13%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.9180"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'-1.29%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.1412"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'1.33%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.04557"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'-7.48%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07088"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'1.18%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.03129"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'-0.04%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.09053"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'-0.95%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.001530"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'-0.48%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.0006164"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'2.05%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.006075"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'-0.74%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'3.459"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'-0.05%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.167"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'0.09%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'2.194"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'1.31%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'-0.19%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'0.84%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'4.091"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'-1.06%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04231"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'-0.22%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.001215"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'0.05%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'-5.86%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'0.21%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'-6.65%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('D40').Value = "'0.03925"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'2.16%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.1115"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'0.02%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.004134"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'6.42%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.002161"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'-10.56%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'-9.29%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.00005165"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'1.23%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'0.09%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.03588"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'-28.17%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'11.56%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.00002100"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'0.09%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0002000"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'0.09%"
$ws.Range('E50').Style = 'Normal'
